# MorganPatrone2006a__M_Stationarygenerator_alpha_zero.xlsx
# "expermits todos no convexos menos el 5to"
#
# Rewrites the numeric / expression values produced by the generator for
# this experiment (leader & follower restriction rows, the modified point,
# vec_bf and vec_BF) while keeping every value stored as TEXT, exactly like
# the original workbook (all these cells are shared-string / text cells,
# never real numbers).

function Set-TextValue($range, $value) {
    # Force the cell to stay text even though the new value looks numeric
    # (mirrors the workbook's original convention of storing these as text).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---------------------------------------------
$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $wsLider.Range("A2") "4.5 - x"
Set-TextValue $wsLider.Range("B2") "-5.0"
Set-TextValue $wsLider.Range("D2") "0.34"
Set-TextValue $wsLider.Range("A3") "-4.5 + x"
Set-TextValue $wsLider.Range("B3") "4.0"
Set-TextValue $wsLider.Range("D3") "0.0"

# --- Restricciones_del_follower -------------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $wsFollower.Range("A2") "2.8 - y"
Set-TextValue $wsFollower.Range("B2") "-3.8"
Set-TextValue $wsFollower.Range("D2") "0.09"
Set-TextValue $wsFollower.Range("E2") "5.0"
Set-TextValue $wsFollower.Range("F2") "7.9"
Set-TextValue $wsFollower.Range("A3") "-2.8 + y"
Set-TextValue $wsFollower.Range("B3") "1.7999999999999998"
Set-TextValue $wsFollower.Range("D3") "0.82"
Set-TextValue $wsFollower.Range("E3") "0"
Set-TextValue $wsFollower.Range("F3") "2.3000000000000003"

# --- Punto_modificado -------------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "4.5"
Set-TextValue $wsPunto.Range("B2") "2.8"

# --- Vector_bf -----------------------------------------------------------
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) resolves case-insensitively (first match wins), so
# name lookups would collide. Use the known sheet indices instead (5 and 6
# respectively - see the workbook's <sheets> order).
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf.Range("A2") "-5.23"

# --- Vector_BF -------------------------------------------------------------
$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF.Range("A2") "1.34"
Set-TextValue $wsBF.Range("A3") "6.0"
